$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pokemon")

# Add new column headers D1:F1 (atk, def, pvMax), matching header style of A1
$ws.Cells.Item(1, 4).Value = "atk"
$ws.Cells.Item(1, 5).Value = "def"
$ws.Cells.Item(1, 6).Value = "pvMax"
$headerStyleSrc = $ws.Range("A1")
$ws.Range("D1:F1").Interior.Color = $headerStyleSrc.Interior.Color

# Pikachu row (row 2): atk=2, def=0, pvMax=20
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 20

# Mewtwo row (row 3): atk=20, def=100, pvMax=200
$ws.Cells.Item(3, 4).Value = 20
$ws.Cells.Item(3, 5).Value = 100
$ws.Cells.Item(3, 6).Value = 200

# Make "pokemon" sheet the active tab and set the selection
[void]$ws.Activate()
[void]$ws.Range("D12").Select()

Write-Host "done"
